$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Messages")

# New row 14: CAN ids assigned
$ws.Cells.Item(14, 1).Value = "ID (Dec)"
$ws.Cells.Item(14, 2).Value = 11
$ws.Cells.Item(14, 3).Value = 12
$ws.Cells.Item(14, 4).Value = 13
$ws.Cells.Item(14, 5).Value = 14
$ws.Cells.Item(14, 7).Value = 26
$ws.Cells.Item(14, 8).Value = 37
$ws.Cells.Item(14, 9).Value = 48
$ws.Cells.Item(14, 10).Value = 59

# Update selection on the Messages sheet
$ws.Range("I16").Select()
